$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The workbook/sheet title tracks the "through" date shown in the report.
$ws.Name = "Through 2022-04-01"

# Capture the current "Total" row (row 5) values before we move anything -
# use Value2 for reliable numeric reads.
$oldB = $ws.Range("B5").Value2
$oldC = $ws.Range("C5").Value2
$oldD = $ws.Range("D5").Value2
$oldE = $ws.Range("E5").Value2
$oldF = $ws.Range("F5").Value2
$oldG = $ws.Range("G5").Value2
$oldH = $ws.Range("H5").Value2
$oldI = $ws.Range("I5").Value2

# New partial-month (April, through 04-01) contributions.
$aprC = 1
$aprE = 1
$aprF = 2
$aprG = 2
$aprH = 1
$aprI = 1

# Relabel row 4 - it was "March (through 03-31)", the month is now complete.
$ws.Range("A4").Value = "March"

# Row 5 becomes the new partial-month April row (sparse - no data in B/D).
$ws.Range("A5").Value = "April (through 04-01)"
$ws.Range("B5").ClearContents()
$ws.Range("C5").Value = $aprC
$ws.Range("D5").ClearContents()
$ws.Range("E5").Value = $aprE
$ws.Range("F5").Value = $aprF
$ws.Range("G5").Value = $aprG
$ws.Range("H5").Value = $aprH
$ws.Range("I5").Value = $aprI

# Row 6 becomes the new "Total" row = old totals + April contributions.
$ws.Range("A6").Value = "Total"
# Copy the bold/bordered label formatting from another label cell onto A6.
$ws.Range("A4").Copy()
$ws.Range("A6").PasteSpecial(-4122)

$ws.Range("B6").Value = $oldB
$ws.Range("C6").Value = $oldC + $aprC
$ws.Range("D6").Value = $oldD
$ws.Range("E6").Value = $oldE + $aprE
$ws.Range("F6").Value = $oldF + $aprF
$ws.Range("G6").Value = $oldG + $aprG
$ws.Range("H6").Value = $oldH + $aprH
$ws.Range("I6").Value = $oldI + $aprI
